# Agregue plantel y se cambio el excel source_informes
#
# Adds a new "Nombre Foto Plantel Club" column (with value
# "Plantel_Rosario_Central" for Diego Martinez) to the "Entrenadores"
# sheet, inserted just before the existing "Nombre Foto Carrera
# Entrenador" column (column I), pushing that column and everything to
# its right one column to the right (I:T -> J:U).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Entrenadores")
$ws.Activate()

# Shift the contents of columns I:T (9..20) one column to the right,
# for every row that currently holds data in that range, working from
# the rightmost column back to column I so that values are never
# clobbered before they are copied.
$lastRow = 3
$firstColToShift = 9   # column I
$lastColToShift = 20   # column T

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = $lastColToShift; $c -ge $firstColToShift; $c--) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r, $c + 1)
        $dstCell.Value2 = $srcCell.Value2
    }
}

# Clear out the (now vacated) old column I value before writing the new
# header/data into it.
$ws.Cells.Item(1, $firstColToShift).Value2 = "Nombre Foto Plantel Club"
$ws.Cells.Item(2, $firstColToShift).Value2 = "Plantel_Rosario_Central"
$ws.Cells.Item(3, $firstColToShift).Value2 = $null

# Update view state: scroll so column B is left-most visible, and move
# the active selection to I3.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("I3").Select()

# Match the print setup used on the other sheet of the workbook.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

Write-Host "Entrenadores sheet updated with Plantel column"
